$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = "VCenter"
$ws.Range("A3").Value = "First Flow"
$ws.Range("A3").Select()
